$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$th = $sm.Theme
Write-Output "theme: $th"
Write-Output "themevariants count: $($th.ThemeVariants.Count)"
